# feat: updated Tommaso Stedile team
# Tommaso Stedile's goalkeeper/team pick (row 21, column B - "Portiere")
# changes from "Nicolas Giordani | FC Savignano" to
# "Federico Leonardi | Sughi ebbasta".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Federico Leonardi | Sughi ebbasta"

# Leave the sheet scrolled back to the top and the cursor on B18, and
# nudge column B a little wider to fit the new (longer) entry - matches
# the view/column-width state captured alongside the data edit.
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$ws.Columns("B").ColumnWidth = 28.8
